# semana 30 de 2025
# Add a new "week 30" column (AG) to the weekly IRA-UCI surveillance sheet,
# mirroring the existing week columns (1..29 in D:AF).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell AG1: the week-number label "30", entered as text (leading
# apostrophe) so it matches the text-typed week headers already in row 1
# (D1="1" ... AF1="29") instead of being auto-coerced to a number.
$ws.Range("AG1").Value = "'30"

# Week-30 counts, one per reporting facility (row). Only rows that already
# carry data through column AF receive a week-30 figure.
$ws.Range("AG2").Value = 0
$ws.Range("AG4").Value = 0
$ws.Range("AG5").Value = 0
$ws.Range("AG6").Value = 1
$ws.Range("AG7").Value = 0
$ws.Range("AG8").Value = 0
$ws.Range("AG9").Value = 0
$ws.Range("AG10").Value = 0
$ws.Range("AG12").Value = 0
$ws.Range("AG14").Value = 0
$ws.Range("AG16").Value = 0
$ws.Range("AG17").Value = 0
$ws.Range("AG22").Value = 0
$ws.Range("AG23").Value = 0
$ws.Range("AG25").Value = 0
$ws.Range("AG28").Value = 53
$ws.Range("AG29").Value = 0
$ws.Range("AG30").Value = 2
$ws.Range("AG31").Value = 0
$ws.Range("AG32").Value = 0
$ws.Range("AG34").Value = 0
$ws.Range("AG35").Value = 3
$ws.Range("AG36").Value = 0
$ws.Range("AG37").Value = 0
$ws.Range("AG38").Value = 0
$ws.Range("AG39").Value = 0
$ws.Range("AG40").Value = 0
$ws.Range("AG41").Value = 0
$ws.Range("AG42").Value = 0
$ws.Range("AG43").Value = 0
$ws.Range("AG44").Value = 0
$ws.Range("AG45").Value = 0
$ws.Range("AG46").Value = 0
$ws.Range("AG47").Value = 0
$ws.Range("AG48").Value = 0
$ws.Range("AG49").Value = 0
$ws.Range("AG50").Value = 0
$ws.Range("AG52").Value = 0
$ws.Range("AG53").Value = 0
$ws.Range("AG54").Value = 0
$ws.Range("AG55").Value = 0
$ws.Range("AG56").Value = 0
$ws.Range("AG57").Value = 0
